$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Price" column (D) values, forcing text storage so numeric-looking
# strings (e.g. trailing zeros, multi-dot groupings) are preserved verbatim.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.229.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.381.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.530"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.380.83"
$ws.Range("D9").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.10"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.809.48"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "59.913.48"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.373.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "322.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "64.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "560.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.494.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0931"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.133"
$ws.Range("D34").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "153.75"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.368"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.07"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0288"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "140.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.587"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0502"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.22"
$ws.Range("D51").Style = "Normal"

# Update "Volume(1h)" column (E) values (plain text, padded with spaces).
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("E3").Value = "  -0.86%  "
$ws.Range("E4").Value = "  +0.69%  "
$ws.Range("E5").Value = "  -1.56%  "
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  +0.67%  "
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("E10").Value = "  -1.61%  "
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("E16").Value = "  -1.62%  "
$ws.Range("E17").Value = "  -1.29%  "
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("E19").Value = "  +11.25%  "
$ws.Range("E21").Value = "  +0.30%  "
$ws.Range("E22").Value = "  +1.14%  "
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("E25").Value = "  -1.12%  "
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("E27").Value = "  -1.61%  "
$ws.Range("E28").Value = "  -4.95%  "
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("E30").Value = "  +2.16%  "
$ws.Range("E31").Value = "  +1.80%  "
$ws.Range("E34").Value = "  +1.27%  "
$ws.Range("E35").Value = "  -0.66%  "
$ws.Range("E36").Value = "  +5.24%  "
$ws.Range("E37").Value = "  +4.10%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("E39").Value = "  -0.75%  "
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("E43").Value = "  +1.87%  "
$ws.Range("E44").Value = "  +0.21%  "
$ws.Range("E45").Value = "  +5.14%  "
$ws.Range("E46").Value = "  +1.51%  "
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("E48").Value = "  +0.99%  "
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("E51").Value = "  -0.27%  "
